$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "SNV gene-level"
$wb.Worksheets.Item(2).Name = "SNV variant-level"
$wb.Worksheets.Item(3).Name = "CNV gene-level"
$wb.Worksheets.Item(4).Name = "Fusion gene-level"
$wb.Worksheets.Item(5).Name = "Fusion fusion-level"
$wb.Worksheets.Item(6).Name = "TPM stats gene-wise z-scores"
$wb.Worksheets.Item(7).Name = "TPM stats group-wise z-scores"
